# "Cập nhật data.xlsx từ công cụ QR"
# The QR tool writes its newest scan result into row 2 (right below the
# header row), pushing all previously-recorded rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data row(s) down and make room for the new entry.
$ws.Rows.Item(2).Insert()

# New QR scan result -> row 2.
$ws.Cells.Item(2, 1).Value  = "ogodx9ovs7q"
$ws.Cells.Item(2, 2).Value  = "kd7046i0"
$ws.Cells.Item(2, 3).Value  = "CTY C"
$ws.Cells.Item(2, 4).Value  = "Madam Thu Bakery, 21C, Võ Văn Tần, Ninh Kiều, Ninh Kiều District, Cần Thơ, 94111, Vietnam"
$ws.Cells.Item(2, 5).Value  = "https://www.google.com/maps/search/?api=1&query=10.032100,105.786400"
$ws.Cells.Item(2, 6).Value  = "2025-08-22T07:01:47.591Z"

# note / phone / branch / cccd / customerCode are unknown for this scan ->
# explicit empty-text cells (matches the blank "" entries used elsewhere in
# this sheet, rather than leaving the cells completely unset).
$ws.Cells.Item(2, 7).Value  = "'"
$ws.Cells.Item(2, 7).ClearFormats()
$ws.Cells.Item(2, 8).Value  = "'"
$ws.Cells.Item(2, 8).ClearFormats()
$ws.Cells.Item(2, 9).Value  = "'"
$ws.Cells.Item(2, 9).ClearFormats()
$ws.Cells.Item(2, 10).Value = "'"
$ws.Cells.Item(2, 10).ClearFormats()
$ws.Cells.Item(2, 11).Value = "'"
$ws.Cells.Item(2, 11).ClearFormats()

$ws.Cells.Item(2, 12).Value = "Nguyễn Văn B"
$ws.Cells.Item(2, 13).Value = "c5fcc4ed3a14b662"
$ws.Cells.Item(2, 14).Value = "d386e590702e53ddeea3640bdb394ee012354447df46f0682febd94210fb1411"
